$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.837.75"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.535.03"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'568.42"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'144.89"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").Value = "2.532.66"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "'27.10"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "2.981.86"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "62.759.98"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "2.527.96"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "'332.70"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "'4.29"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "'65.01"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "'1.57"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "'7.21"
$ws.Range("E30").Value = "  +5.19%  "
$ws.Range("D31").Value = "0.0₃0801"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D33").Value = "'176.90"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "'394.90"
$ws.Range("E35").Value = "  -5.10%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'39.29"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("D43").Value = "'150.36"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "'0.0526"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "'0.0959"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'17.95"
$ws.Range("E50").Value = "  -5.73%  "
$ws.Range("D51").Value = "'11.30"
$ws.Range("E51").Value = "  +0.38%  "
